$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "27.834.52"
Set-TextCell "E2" "  +3.29%  "
Set-TextCell "D3" "1.868.24"
Set-TextCell "E3" "  +2.90%  "
Set-TextCell "D4" "1.041"
Set-TextCell "E4" "  +3.52%  "
Set-TextCell "D5" "324.54"
Set-TextCell "E5" "  +4.00%  "
Set-TextCell "D6" "1.037"
Set-TextCell "D7" "0.4431"
Set-TextCell "E7" "  +3.16%  "
Set-TextCell "D8" "0.3804"
Set-TextCell "E8" "  +3.65%  "
Set-TextCell "D9" "0.07479"
Set-TextCell "E9" "  +3.37%  "
Set-TextCell "E10" "  +2.61%  "
Set-TextCell "D11" "21.80"
Set-TextCell "D12" "1.884.20"
Set-TextCell "E12" "  -11.44%  "
Set-TextCell "D13" "5.567"
Set-TextCell "E13" "  +2.86%  "
Set-TextCell "D14" "6.772"
Set-TextCell "E14" "  +2.51%  "
Set-TextCell "D15" "0.07232"
Set-TextCell "E15" "  +4.02%  "
Set-TextCell "D16" "83.88"
Set-TextCell "E16" "  +3.48%  "
Set-TextCell "E17" "  +3.23%  "
Set-TextCell "D18" "0.000009171"
Set-TextCell "E18" "  +3.34%  "
Set-TextCell "D19" "1.038"
Set-TextCell "D20" "15.58"
Set-TextCell "E20" "  +2.12%  "
Set-TextCell "D21" "27.860.01"
Set-TextCell "E21" "  +3.20%  "
Set-TextCell "D22" "5.334"
Set-TextCell "E22" "  +2.89%  "
Set-TextCell "D23" "11.37"
Set-TextCell "E23" "  +3.58%  "
Set-TextCell "D24" "1.983"
Set-TextCell "E24" "  +5.34%  "
Set-TextCell "D25" "158.79"
Set-TextCell "E25" "  +3.12%  "
Set-TextCell "D26" "18.91"
Set-TextCell "E26" "  +3.08%  "
Set-TextCell "B27" "InternetComputer(DFINITY)"
Set-TextCell "C27" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D27" "5.348"
Set-TextCell "E27" "  +2.45%  "
Set-TextCell "B28" "LidoDAOToken"
Set-TextCell "C28" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D28" "1.992"
Set-TextCell "E28" "  +4.87%  "
Set-TextCell "D29" "117.90"
Set-TextCell "E29" "  +2.87%  "
Set-TextCell "D30" "0.09092"
Set-TextCell "E30" "  +1.64%  "
Set-TextCell "D31" "0.7790"
Set-TextCell "E31" "  +4.27%  "
Set-TextCell "D32" "3.115"
Set-TextCell "E32" "  +10.86%  "
Set-TextCell "D33" "1.219"
Set-TextCell "E33" "  +2.51%  "
Set-TextCell "D34" "4.584"
Set-TextCell "E34" "  +3.77%  "
Set-TextCell "D35" "1.039"
Set-TextCell "E35" "  +3.23%  "
Set-TextCell "E36" "  +2.52%  "
Set-TextCell "D37" "0.01999"
Set-TextCell "E37" "  +4.02%  "
Set-TextCell "D38" "0.05358"
Set-TextCell "E38" "  +2.71%  "
Set-TextCell "E39" "  +4.44%  "
Set-TextCell "D40" "0.5209"
Set-TextCell "E40" "  +2.15%  "
Set-TextCell "E41" "  +2.64%  "
Set-TextCell "D42" "6.929"
Set-TextCell "E42" "  +7.13%  "
Set-TextCell "D43" "8.696"
Set-TextCell "E43" "  +4.71%  "
Set-TextCell "D44" "109.81"
Set-TextCell "E44" "  +2.83%  "
Set-TextCell "D45" "10.70"
Set-TextCell "E45" "  +2.63%  "
Set-TextCell "D46" "1.732"
Set-TextCell "E46" "  +5.25%  "
Set-TextCell "D47" "0.4720"
Set-TextCell "E47" "  +3.38%  "
Set-TextCell "D48" "0.06463"
Set-TextCell "E48" "  +3.93%  "
Set-TextCell "D49" "1.907"
Set-TextCell "D50" "39.93"
Set-TextCell "E50" "  +4.10%  "
Set-TextCell "D51" "64.81"
Set-TextCell "E51" "  +2.78%  "
